$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.871.72'
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").Value = '1.731.48'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D4").Value = '''0.9973'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''242.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '''0.9978'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '''0.4897'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").Value = '''0.2600'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("D9").Value = '''0.06220'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").Value = '1.735.57'
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("D11").Value = '''16.05'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.30%  '
$ws.Range("D12").Value = '''0.06892'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.37%  '
$ws.Range("D13").Value = '''0.6106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '''4.507'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").Value = '''77.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '''0.9985'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '26.636.30'
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").Value = '''0.9973'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '''0.000007182'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '''11.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("D21").Value = '1.958.74'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '''4.437'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '''8.572'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '''5.110'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").Value = '''138.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").Value = '''15.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = '''1.776'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.89%  '
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''106.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''1.380'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").Value = '''3.926'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").Value = '''0.07982'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '''3.678'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '''0.04522'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '''0.9970'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.604'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''1.010'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.6231'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''0.9441'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.28%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''2.048'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.435'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '''0.9979'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.01507'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.628'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.34%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''99.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.3860'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '''6.929'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1162'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.05389'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.50%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.873'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''30.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.240'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.11%  '
